$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 5.2
$ws.Range("J2").Value = 2.27
$ws.Range("K2").Value = 2.07
$ws.Range("L2").Value = 5.1
$ws.Range("M2").Value = 9.35
$ws.Range("N2").Value = 1.03
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 3.1
$ws.Range("Q2").Value = 1.87
$ws.Range("R2").Value = 1.83
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.52
$ws.Range("U2").Value = 1.78
$ws.Range("V2").Value = 1.83
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 13.5
$ws.Range("AC2").Value = 9.75
$ws.Range("AD2").Value = 6.8
$ws.Range("AE2").Value = 15
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 35
$ws.Range("AJ2").Value = 16
$ws.Range("AK2").Value = 110
$ws.Range("AL2").Value = 50
$ws.Range("AN2").Value = 3.45
$ws.Range("AO2").Value = 8.5
$ws.Range("AQ2").Value = 29
$ws.Range("AT2").Value = 2.47
$ws.Range("AU2").Value = 7.2
$ws.Range("AW2").Value = 6.6
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 32
$ws.Range("BA2").Value = 200
